# Auto-generated edit script: applies the scheduled-runner profit-recalculation update
# to the Sheets workbook. For each touched row, H..N (currentAveragePrice.. LeveProfitHQ)
# are refreshed to the latest recalculated values. A few rows gain or lose a cell
# (LeveProfitNQ / LeveProfitHQ) when the HQ/NQ branch of the calc toggles on/off.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ALC_updates = @(
  @{Row=70; Col=8; Val=6482.619},
  @{Row=70; Col=9; Val=3666.1667},
  @{Row=70; Col=10; Val=7609.2},
  @{Row=70; Col=11; Val=10998.5001},
  @{Row=70; Col=12; Val=22827.6},
  @{Row=70; Col=13; Val=-10728.5001},
  @{Row=70; Col=14; Val=-23367.6},
  @{Row=73; Col=8; Val=6482.619},
  @{Row=73; Col=9; Val=3666.1667},
  @{Row=73; Col=10; Val=7609.2},
  @{Row=73; Col=11; Val=10998.5001},
  @{Row=73; Col=12; Val=22827.6},
  @{Row=73; Col=13; Val=-10062.5001},
  @{Row=73; Col=14; Val=-24699.6},
  @{Row=80; Col=8; Val=3124.8975},
  @{Row=80; Col=9; Val=1672.6842},
  @{Row=80; Col=11; Val=5018.0526},
  @{Row=80; Col=13; Val=-4020.0526},
  @{Row=82; Col=8; Val=346742.16},
  @{Row=82; Col=9; Val=346742.16},
  @{Row=82; Col=11; Val=1040226.48},
  @{Row=82; Col=13; Val=-1039820.48},
  @{Row=83; Col=8; Val=3124.8975},
  @{Row=83; Col=9; Val=1672.6842},
  @{Row=83; Col=11; Val=15054.1578},
  @{Row=83; Col=13; Val=-10062.1578},
  @{Row=85; Col=8; Val=346742.16},
  @{Row=85; Col=9; Val=346742.16},
  @{Row=85; Col=11; Val=1040226.48},
  @{Row=85; Col=13; Val=-1038822.48},
  @{Row=92; Col=8; Val=1247.5},
  @{Row=92; Col=9; Val=1247.5},
  @{Row=92; Col=11; Val=1247.5},
  @{Row=92; Col=13; Val=0.5},
  @{Row=125; Col=8; Val=2433},
  @{Row=125; Col=10; Val=3024.5},
  @{Row=125; Col=12; Val=27220.5},
  @{Row=125; Col=14; Val=-32140.5},
  @{Row=127; Col=8; Val=11543.533},
  @{Row=127; Col=9; Val=2295.75},
  @{Row=127; Col=10; Val=22112.428},
  @{Row=127; Col=11; Val=6887.25},
  @{Row=127; Col=12; Val=66337.284},
  @{Row=127; Col=13; Val=-1927.25},
  @{Row=127; Col=14; Val=-76257.284},
  @{Row=131; Col=8; Val=8297.259},
  @{Row=131; Col=10; Val=7936.4614},
  @{Row=131; Col=12; Val=23809.3842},
  @{Row=131; Col=14; Val=-33889.3842},
  @{Row=132; Col=8; Val=19362.133},
  @{Row=132; Col=9; Val=1301.8182},
  @{Row=132; Col=11; Val=3905.4546},
  @{Row=132; Col=13; Val=-1375.4546},
  @{Row=135; Col=8; Val=4612.45},
  @{Row=135; Col=9; Val=4646.7144},
  @{Row=135; Col=10; Val=4532.5},
  @{Row=135; Col=11; Val=41820.4296},
  @{Row=135; Col=12; Val=40792.5},
  @{Row=135; Col=13; Val=-39285.4296},
  @{Row=135; Col=14; Val=-45862.5},
  @{Row=137; Col=8; Val=6461461.5},
  @{Row=137; Col=10; Val=2327.3},
  @{Row=137; Col=12; Val=6981.900000000001},
  @{Row=137; Col=14; Val=-12081.9},
  @{Row=138; Col=8; Val=4853.65},
  @{Row=138; Col=9; Val=4462.6665},
  @{Row=138; Col=10; Val=5021.2144},
  @{Row=138; Col=11; Val=13387.9995},
  @{Row=138; Col=12; Val=15063.6432},
  @{Row=138; Col=13; Val=-8247.999500000002},
  @{Row=138; Col=14; Val=-25343.6432},
  @{Row=141; Col=8; Val=6389.39},
  @{Row=141; Col=9; Val=3999.0303},
  @{Row=141; Col=11; Val=11997.0909},
  @{Row=141; Col=13; Val=-6817.090899999999}
)
foreach ($u in $ALC_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ARM_updates = @(
  @{Row=74; Col=8; Val=3703.7693},
  @{Row=74; Col=9; Val=3703.7693},
  @{Row=74; Col=10; Val=0},
  @{Row=74; Col=11; Val=3703.7693},
  @{Row=74; Col=12; Val=0},
  @{Row=74; Col=13; Val=-2829.7693},
  @{Row=77; Col=8; Val=3703.7693},
  @{Row=77; Col=9; Val=3703.7693},
  @{Row=77; Col=10; Val=0},
  @{Row=77; Col=11; Val=18518.8465},
  @{Row=77; Col=12; Val=0},
  @{Row=77; Col=13; Val=-14150.8465}
)
foreach ($u in $ARM_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

$ARM_clears = @(
  @{Row=74; Col=14},
  @{Row=77; Col=14}
)
foreach ($u in $ARM_clears) { $ws.Cells.Item($u.Row, $u.Col).ClearContents() }

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

$BSM_updates = @(
  @{Row=86; Col=8; Val=49776.094},
  @{Row=86; Col=10; Val=2377.111},
  @{Row=86; Col=12; Val=2377.111},
  @{Row=86; Col=14; Val=-4623.111},
  @{Row=89; Col=8; Val=49776.094},
  @{Row=89; Col=10; Val=2377.111},
  @{Row=89; Col=12; Val=11885.555},
  @{Row=89; Col=14; Val=-23117.555},
  @{Row=94; Col=8; Val=2693.375},
  @{Row=94; Col=10; Val=5283.75},
  @{Row=94; Col=12; Val=5283.75},
  @{Row=94; Col=14; Val=-6185.75},
  @{Row=107; Col=8; Val=4801.0454},
  @{Row=107; Col=10; Val=6600},
  @{Row=107; Col=12; Val=6600},
  @{Row=107; Col=14; Val=-10440},
  @{Row=134; Col=8; Val=3636.8708},
  @{Row=134; Col=9; Val=2165.9622},
  @{Row=134; Col=11; Val=6497.8866},
  @{Row=134; Col=13; Val=-3962.8866}
)
foreach ($u in $BSM_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

$CRP_updates = @(
  @{Row=31; Col=8; Val=4801.0713},
  @{Row=31; Col=9; Val=5446.8887},
  @{Row=31; Col=10; Val=3638.6},
  @{Row=31; Col=11; Val=5446.8887},
  @{Row=31; Col=12; Val=3638.6},
  @{Row=31; Col=13; Val=-5151.8887},
  @{Row=31; Col=14; Val=-4228.6},
  @{Row=34; Col=8; Val=4801.0713},
  @{Row=34; Col=9; Val=5446.8887},
  @{Row=34; Col=10; Val=3638.6},
  @{Row=34; Col=11; Val=5446.8887},
  @{Row=34; Col=12; Val=3638.6},
  @{Row=34; Col=13; Val=-5244.8887},
  @{Row=34; Col=14; Val=-4042.6},
  @{Row=97; Col=8; Val=42048.25},
  @{Row=97; Col=9; Val=20000},
  @{Row=97; Col=10; Val=49397.668},
  @{Row=97; Col=11; Val=20000},
  @{Row=97; Col=12; Val=49397.668},
  @{Row=97; Col=13; Val=-19009},
  @{Row=97; Col=14; Val=-51379.668},
  @{Row=132; Col=8; Val=851.25},
  @{Row=132; Col=9; Val=851.25},
  @{Row=132; Col=11; Val=2553.75},
  @{Row=132; Col=13; Val=-23.75}
)
foreach ($u in $CRP_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

$CUL_updates = @(
  @{Row=68; Col=8; Val=2831},
  @{Row=68; Col=10; Val=3073},
  @{Row=68; Col=12; Val=9219},
  @{Row=68; Col=14; Val=-10841},
  @{Row=71; Col=8; Val=2831},
  @{Row=71; Col=10; Val=3073},
  @{Row=71; Col=12; Val=27657},
  @{Row=71; Col=14; Val=-35769},
  @{Row=130; Col=8; Val=1168666.6},
  @{Row=130; Col=9; Val=1751500},
  @{Row=130; Col=10; Val=3000},
  @{Row=130; Col=11; Val=5254500},
  @{Row=130; Col=12; Val=9000},
  @{Row=130; Col=13; Val=-5249480},
  @{Row=130; Col=14; Val=-19040},
  @{Row=137; Col=8; Val=2482},
  @{Row=137; Col=9; Val=1465},
  @{Row=137; Col=11; Val=4395},
  @{Row=137; Col=13; Val=705}
)
foreach ($u in $CUL_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

$GSM_updates = @(
  @{Row=64; Col=8; Val=40001},
  @{Row=64; Col=10; Val=40001},
  @{Row=64; Col=12; Val=40001},
  @{Row=64; Col=14; Val=-40497},
  @{Row=67; Col=8; Val=40001},
  @{Row=67; Col=10; Val=40001},
  @{Row=67; Col=12; Val=40001},
  @{Row=67; Col=14; Val=-41717},
  @{Row=70; Col=8; Val=70283.64999999999},
  @{Row=70; Col=9; Val=127502.22},
  @{Row=70; Col=10; Val=5912.75},
  @{Row=70; Col=11; Val=127502.22},
  @{Row=70; Col=12; Val=5912.75},
  @{Row=70; Col=13; Val=-127232.22},
  @{Row=70; Col=14; Val=-6452.75},
  @{Row=73; Col=8; Val=70283.64999999999},
  @{Row=73; Col=9; Val=127502.22},
  @{Row=73; Col=10; Val=5912.75},
  @{Row=73; Col=11; Val=127502.22},
  @{Row=73; Col=12; Val=5912.75},
  @{Row=73; Col=13; Val=-126566.22},
  @{Row=73; Col=14; Val=-7784.75},
  @{Row=113; Col=8; Val=11760.615},
  @{Row=113; Col=10; Val=13611},
  @{Row=113; Col=12; Val=13611},
  @{Row=113; Col=14; Val=-17951},
  @{Row=122; Col=8; Val=4849.5713},
  @{Row=122; Col=9; Val=4407.8335},
  @{Row=122; Col=10; Val=7500},
  @{Row=122; Col=11; Val=13223.5005},
  @{Row=122; Col=12; Val=22500},
  @{Row=122; Col=13; Val=-10773.5005},
  @{Row=122; Col=14; Val=-27400},
  @{Row=126; Col=8; Val=3994.6},
  @{Row=126; Col=9; Val=3993.8572},
  @{Row=126; Col=10; Val=3996.3333},
  @{Row=126; Col=11; Val=11981.5716},
  @{Row=126; Col=12; Val=11988.9999},
  @{Row=126; Col=13; Val=-9511.571599999999},
  @{Row=126; Col=14; Val=-16928.9999}
)
foreach ($u in $GSM_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

$LTW_updates = @(
  @{Row=7; Col=8; Val=6428.857},
  @{Row=7; Col=9; Val=6417.1665},
  @{Row=7; Col=10; Val=6499},
  @{Row=7; Col=11; Val=6417.1665},
  @{Row=7; Col=12; Val=6499},
  @{Row=7; Col=13; Val=-6305.1665},
  @{Row=7; Col=14; Val=-6723},
  @{Row=22; Col=8; Val=0},
  @{Row=22; Col=10; Val=0},
  @{Row=22; Col=12; Val=0},
  @{Row=27; Col=8; Val=0},
  @{Row=27; Col=10; Val=0},
  @{Row=27; Col=12; Val=0},
  @{Row=40; Col=8; Val=1995.3572},
  @{Row=40; Col=9; Val=1995.3572},
  @{Row=40; Col=10; Val=0},
  @{Row=40; Col=11; Val=1995.3572},
  @{Row=40; Col=12; Val=0},
  @{Row=40; Col=13; Val=-1859.3572},
  @{Row=123; Col=8; Val=28357.143},
  @{Row=123; Col=9; Val=10000},
  @{Row=123; Col=10; Val=29769.23},
  @{Row=123; Col=11; Val=10000},
  @{Row=123; Col=12; Val=29769.23},
  @{Row=123; Col=13; Val=-5100},
  @{Row=123; Col=14; Val=-39569.23},
  @{Row=126; Col=8; Val=6428.857},
  @{Row=126; Col=9; Val=6417.1665},
  @{Row=126; Col=10; Val=6499},
  @{Row=126; Col=11; Val=19251.4995},
  @{Row=126; Col=12; Val=19497},
  @{Row=126; Col=13; Val=-16781.4995},
  @{Row=126; Col=14; Val=-24437}
)
foreach ($u in $LTW_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

$LTW_clears = @(
  @{Row=22; Col=14},
  @{Row=27; Col=14},
  @{Row=40; Col=14}
)
foreach ($u in $LTW_clears) { $ws.Cells.Item($u.Row, $u.Col).ClearContents() }

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

$WVR_updates = @(
  @{Row=40; Col=8; Val=21497.5},
  @{Row=40; Col=10; Val=21497.5},
  @{Row=40; Col=12; Val=21497.5},
  @{Row=40; Col=14; Val=-21795.5},
  @{Row=74; Col=8; Val=59990},
  @{Row=74; Col=10; Val=59990},
  @{Row=74; Col=12; Val=59990},
  @{Row=74; Col=14; Val=-61862},
  @{Row=77; Col=8; Val=59990},
  @{Row=77; Col=10; Val=59990},
  @{Row=77; Col=12; Val=179970},
  @{Row=77; Col=14; Val=-189330},
  @{Row=94; Col=8; Val=25000},
  @{Row=94; Col=10; Val=25000},
  @{Row=94; Col=12; Val=25000},
  @{Row=94; Col=14; Val=-26802},
  @{Row=115; Col=8; Val=25992.309},
  @{Row=115; Col=10; Val=25992.309},
  @{Row=115; Col=12; Val=25992.309},
  @{Row=115; Col=14; Val=-29126.309},
  @{Row=132; Col=8; Val=2109.7},
  @{Row=132; Col=9; Val=2109.7},
  @{Row=132; Col=11; Val=6329.099999999999},
  @{Row=132; Col=13; Val=-3799.099999999999}
)
foreach ($u in $WVR_updates) { $ws.Cells.Item($u.Row, $u.Col).Value = $u.Val }

Write-Host "Applied $(248) cell updates and $(5) clears"
